# Update the 'F' column (想去人数 / interested-count) figures to match
# the freshly generated gh-pages data snapshot (commit 456a3b4).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 847
$ws.Range("F5").Value = 524
$ws.Range("F6").Value = 1172
$ws.Range("F10").Value = 128
$ws.Range("F11").Value = 1237
$ws.Range("F14").Value = 925
$ws.Range("F15").Value = 900
$ws.Range("F16").Value = 189
$ws.Range("F20").Value = 842
$ws.Range("F21").Value = 1764
$ws.Range("F22").Value = 3299
$ws.Range("F23").Value = 974
$ws.Range("F24").Value = 95
$ws.Range("F25").Value = 2364
$ws.Range("F27").Value = 25
$ws.Range("F28").Value = 3251
$ws.Range("F29").Value = 679
$ws.Range("F32").Value = 2012
$ws.Range("F34").Value = 766
$ws.Range("F35").Value = 154
$ws.Range("F36").Value = 147
$ws.Range("F37").Value = 98
$ws.Range("F39").Value = 1172
$ws.Range("F40").Value = 1842
$ws.Range("F41").Value = 436
$ws.Range("F44").Value = 220
$ws.Range("F46").Value = 202

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 145

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 165

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 847
$ws.Range("F5").Value = 1172
$ws.Range("F7").Value = 1237
$ws.Range("F8").Value = 925
$ws.Range("F9").Value = 900
$ws.Range("F10").Value = 145
$ws.Range("F16").Value = 842
$ws.Range("F17").Value = 1765
$ws.Range("F18").Value = 3299
$ws.Range("F19").Value = 974
$ws.Range("F20").Value = 95
$ws.Range("F21").Value = 2364
$ws.Range("F22").Value = 25
$ws.Range("F23").Value = 3251
$ws.Range("F24").Value = 679
$ws.Range("F28").Value = 2012
$ws.Range("F34").Value = 766
$ws.Range("F35").Value = 154
$ws.Range("F36").Value = 147
$ws.Range("F37").Value = 98
$ws.Range("F41").Value = 1172
$ws.Range("F42").Value = 1842
$ws.Range("F45").Value = 436
$ws.Range("F47").Value = 220
$ws.Range("F49").Value = 202
